$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename Sheet1 -> test.csv (fix csv/xls unicode support)
$ws.Name = "test.csv"

# Insert two new columns: B (Sparse Data) and D (Unicode Test)
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(4).Insert()

# Header row
$ws.Range("B1").Value = "Sparse Data"
$ws.Range("D1").Value = "Unicode Test"

# Data row 2 - unicode values
$ws.Range("B2").Value = "Iñtërnâtiônàližætiøn"
$ws.Range("D2").Value = "Ādam"

# The column insert bled A2's date style into B2; clear it back to Normal
# so the new string cell in B2 carries no number-format style, matching
# the sparse nature of the new columns.
$ws.Range("B2").Style = "Normal"

# Rows 3 and 4 stay empty for the new sparse columns, but they keep the
# same number-format style as column A on those rows (this is what Excel's
# column-insert does naturally for the cells it creates).
$ws.Range("A3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("A4").Copy()
$ws.Range("D4").PasteSpecial(-4122)

# Give the new column B/D the same on-disk width as column A (16).
# Excel's ColumnWidth property is in character units and gets offset by the
# font's default-character-width padding when serialized to the raw <col>
# width attribute, so we back the 16 target out through that same offset.
$ws.Columns.Item(2).ColumnWidth = 15.166666666666666
$ws.Columns.Item(4).ColumnWidth = 15.166666666666666

# Move the active selection to D2, as in the edited workbook
$ws.Range("D2").Select()

Write-Output "done"
